$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value  = -7.7135
$ws.Range("D6").Value  = -7.891699999999997
$ws.Range("C7").Value  = -11.19769999999999
$ws.Range("D7").Value  = -7.755999999999996
$ws.Range("A8").Value  = -20.95130000000001
$ws.Range("D8").Value  = -7.993100000000002
$ws.Range("D9").Value  = -8.176000000000005
$ws.Range("A10").Value = -20.49489999999998
$ws.Range("D10").Value = -6.523699999999998
$ws.Range("A12").Value = -22.56430000000004
$ws.Range("D12").Value = -8.329300000000007
$ws.Range("B13").Value = 6.183199999999999
$ws.Range("A18").Value = -22.52230000000004
$ws.Range("C20").Value = -15.01009999999999
$ws.Range("A25").Value = -22.25800000000003
